# Update the 'K' column (column G) values for rows 2-29 to reflect the
# regenerated save_data (K computed from Strike# instead of the old value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 7
    4  = 8
    5  = 6
    6  = 6
    7  = 8
    8  = 10
    9  = 5
    10 = 5
    11 = 5
    12 = 8
    13 = 4
    14 = 9
    15 = 6
    16 = 5
    17 = 7
    18 = 5
    19 = 3
    20 = 8
    21 = 6
    22 = 12
    23 = 7
    24 = 5
    25 = 5
    26 = 5
    27 = 6
    28 = 1
    29 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
